# feat: add 2022-Q4 data
#
# The workbook tracks quarterly holdings. A new quarter ("2022-Q4") is
# inserted right after the "总计" (totals) summary sheet. The sheet that
# used to hold the most-recent quarter ("2022-Q3") is refreshed in place
# with the new "2022-Q4" numbers and renamed, while a fresh copy of it
# (still carrying the old "2022-Q3" numbers) is placed right after it so
# the quarterly archive is preserved. The summary sheet gains a new row
# for "2022-Q4" and the row that used to be missing ("2021-Q2") is
# appended at the bottom.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate the current "2022-Q3" sheet so its existing figures are
#    preserved under the same name, then refresh the original sheet
#    with the new "2022-Q4" figures and rename it.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $q3)
$q3archive = $wb.Worksheets.Item("2022-Q3 (2)")

$q4 = $q3
$q4.Name = "2022-Q4"
$q3archive.Name = "2022-Q3"

$q4.Range("D2").Value = "'2.81"
$q4.Range("E2").Value = "'93.63"
$q4.Range("F2").Value = "'3.99"
$q4.Range("G2").Value = "'0.1121"
$q4.Range("H2").Value = 5

$q4.Range("D3").Value = "'0.64"
$q4.Range("E3").Value = "'93.56"
$q4.Range("F3").Value = "'7.01"
$q4.Range("G3").Value = "'0.0449"

$q4.Range("D4").Value = "'0.43"
$q4.Range("E4").Value = "'92.90"
$q4.Range("F4").Value = "'3.64"
$q4.Range("G4").Value = "'0.0157"
$q4.Range("H4").Value = 2

$q4.Range("D5").Value = "'0.36"
$q4.Range("E5").Value = "'93.63"
$q4.Range("F5").Value = "'3.99"
$q4.Range("G5").Value = "'0.0144"
$q4.Range("H5").Value = 5

# ---------------------------------------------------------------------
# 2) Update the "总计" (totals) summary sheet: shift the quarter labels
#    down one row, add the new "2022-Q4" entry at the top of the data,
#    and append the previously-missing "2021-Q2" row at the bottom.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.19

$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 0.16

$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 5
$total.Range("D4").Value = 0.32

$total.Range("B5").Value = "2021-Q3"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.24

$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)
$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q2"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.27

# ---------------------------------------------------------------------
# 3) Restore the originally-selected tab ("2021-Q2", the last sheet).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
